$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for the new columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the style from an existing header cell (AC1) to the new header cells
$headerStyle = $ws.Range("AC1")
$headerStyle.Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in team record data for all data rows (2 through 51)
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 30).Value = 81  # AD
    $ws.Cells.Item($r, 31).Value = 81  # AE
    $ws.Cells.Item($r, 32).Value = 0   # AF
}
